$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "-3.13, 0.892"
$ws.Range("D7").Value = "-0.819, 1.87"
$ws.Range("D8").Value = "-1.785, 0.918"
$ws.Range("C9").Value = 0.358
$ws.Range("D9").Value = "-0.817, 0.586"
$ws.Range("D10").Value = "-0.311, 0.289"
